# Hid the first two slides
$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$s1.SlideShowTransition.Hidden = 1

$s2 = $p.Slides.Item(2)
$s2.SlideShowTransition.Hidden = 1
